# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.546.66"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.108.00"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.17"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5227"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4536"
$ws.Range("E8").Value = "  +3.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.41"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09022"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.170"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.58"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.110.13"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.826"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.110"
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001170"
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.11"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.011"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06683"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.33"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.228"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.602.92"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.81"
$ws.Range("E24").Value = "  +4.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.353.97"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.36"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.48"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.509"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.36"
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.214"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.636"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.335"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.41"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.880"
$ws.Range("E37").Value = "  +8.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02614"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2309"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6832"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.254"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6430"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.07"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.304"
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.673"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000356"
$ws.Range("E48").Value = "  +19.49%  "
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3386"
$ws.Range("E50").Value = "  +12.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.01"
$ws.Range("E51").Value = "  +1.55%  "
